# Auto commit at 2026-01-13  7:52:33.90
# Updates the monthly metrics figures (Metrics sheet) and refreshes the
# "today" sheet's running totals, which derive from Metrics via formulas.
# The per-category daily increments on "today" (B3:B6) are cleared because
# they are not part of today's running total anymore.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refreshed monthly figures ---------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 229710.65
$metrics.Range("B3").Value  = 171400.75
$metrics.Range("B4").Value  = 58644.900000000009
$metrics.Range("B5").Value  = 9309
$metrics.Range("B6").Value  = 5865581.379999999
$metrics.Range("B7").Value  = 4942118.38
$metrics.Range("B8").Value  = 1722736.72
$metrics.Range("B9").Value  = 229586
$metrics.Range("B10").Value = 34330962.369999997
$metrics.Range("B11").Value = 32217393.539999999
$metrics.Range("B12").Value = 12004458.76
$metrics.Range("B13").Value = 1327216

# --- "today" sheet: clear yesterday's incremental figures ---------------
$today = $wb.Worksheets.Item("today")

$today.Range("B3").ClearContents()
$today.Range("B4").ClearContents()
$today.Range("B5").ClearContents()
$today.Range("B6").ClearContents()

# --- Restore the selections recorded in each sheet's view ----------------
# ("today" was left selected first so that "Metrics" ends up as the
# workbook's active tab, matching the saved file.)
[void]$today.Activate()
[void]$today.Range("E9").Select()

[void]$metrics.Activate()
[void]$metrics.Range("C32").Select()
